$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet to reflect the unified "DataNode" concept
$ws.Name = "DataNode"

# Restore the author's last selection before saving
$ws.Range("D36").Select()
